# Fix: 3T segmentation regions were swapped w.r.t. the 1.5T segmentation.
# Swap the temperature (D) and uncertainty (E) values between each pair of
# consecutive rows (region_id 1 and region_id 2) for every "run".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows run from 2 to 17, in pairs: (2,3), (4,5), (6,7), ..., (16,17)
for ($row = 2; $row -le 16; $row += 2) {
    $row2 = $row + 1

    $dValue1 = $ws.Cells.Item($row, 4).Value2
    $eValue1 = $ws.Cells.Item($row, 5).Value2

    $dValue2 = $ws.Cells.Item($row2, 4).Value2
    $eValue2 = $ws.Cells.Item($row2, 5).Value2

    $ws.Cells.Item($row, 4).Value2 = $dValue2
    $ws.Cells.Item($row, 5).Value2 = $eValue2

    $ws.Cells.Item($row2, 4).Value2 = $dValue1
    $ws.Cells.Item($row2, 5).Value2 = $eValue1
}

$wb.Save()
